# Update cryptocurrency price/volume data per latest refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "64.889.74"
$ws.Range("E2").Value = "  -0.54%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.563.28"
$ws.Range("E3").Value = "  +2.46%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.03%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.13"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.86%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.54%  "

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "3.562.26"
$ws.Range("E7").Value = "  +2.45%  "

# Row 8: USDC
$ws.Range("E8").Value = "  -0.01%  "

# Row 9: XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.57%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  +0.53%  "

# Row 11: Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.96"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.89%  "

# Row 12: Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.384"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.55%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.171.90"
$ws.Range("E13").Value = "  +2.47%  "

# Row 14: ShibaInu
$ws.Range("E14").Value = "  +0.19%  "

# Row 15: WrappedEther
$ws.Range("D15").Value = "3.565.14"
$ws.Range("E15").Value = "  +2.19%  "

# Row 16: Avalanche
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.00"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.87%  "

# Row 17: TRON
$ws.Range("E17").Value = "  +0.46%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "65.042.37"
$ws.Range("E18").Value = "  -0.15%  "

# Row 19: Uniswap
$ws.Range("E19").Value = "  +3.20%  "

# Row 20: Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.35"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.49%  "

# Row 21: Polkadot
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.88%  "

# Row 22: BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "388.61"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.09%  "

# Row 23: Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.579"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.64%  "

# Row 24: WrappedeETH
$ws.Range("D24").Value = "3.711.03"
$ws.Range("E24").Value = "  +2.52%  "

# Row 25: Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.92%  "

# Row 26: Dai
$ws.Range("E26").Value = "  +0.06%  "

# Row 27: PEPE
$ws.Range("E27").Value = "  +5.46%  "

# Row 28: RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.72"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.93%  "

# Row 29: Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.10%  "

# Row 30: PancakeSwap
$ws.Range("E30").Value = "  +3.24%  "

# Row 31: InternetComputer(DFINITY)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.42"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.93%  "

# Row 32: Fetch.AI
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.48"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +24.25%  "

# Row 33: RenzoRestakedETH
$ws.Range("D33").Value = "3.566.94"
$ws.Range("E33").Value = "  +1.96%  "

# Row 34: EthereumClassic
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.02"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.01%  "

# Row 35: USDe
$ws.Range("E35").Value = "  -0.01%  "

# Row 36: Kaspa
$ws.Range("E36").Value = "  +0.36%  "

# Row 37: Aptos -> Monero
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "169.18"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.57%  "

# Row 38: Monero -> Aptos
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.92"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.73%  "

# Row 39: ImmutableX
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.54"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.92%  "

# Row 40: NEARProtocol
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.98"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +5.67%  "

# Row 41: Hedera
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0807"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +3.57%  "

# Row 42: EnergySwap
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "27.04"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.59%  "

# Row 43: Mantle
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.825"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.83%  "

# Row 44: OKB
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.71"
$ws.Range("D44").ClearFormats()

# Row 46: Filecoin
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.46"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.53%  "

# Row 47: ONDO
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.20"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.13%  "

# Row 48: Stacks
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.64"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.37%  "

# Row 49: Maker
$ws.Range("D49").Value = "2.478.94"
$ws.Range("E49").Value = "  +11.59%  "

# Row 50: Cosmos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.89"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.09%  "

# Row 51: LidoDAOToken -> SuiNetwork
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.867"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.74%  "
